$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 from "s" to "asd"
$ws.Cells.Item(2, 2).Value = "asd"

# Update C2 from numeric 1 to the text string "2" (force text type, then
# restore the default "Normal" style so no stray number-format style sticks)
$c2 = $ws.Cells.Item(2, 3)
$c2.NumberFormat = "@"
$c2.Value = "2"
$c2.Style = "Normal"

# Delete row 3 entirely (the Expense/asd/1 row)
$ws.Rows("3:3").Delete()
